# Changes of 25th May 2022
# Updates the PackageTrackNum (column C) values for rows 2-22, and for the
# rows where ShipmentTrackNum (column D) mirrors the PackageTrackNum value
# (rows 5-7 and 13-17), updates column D to match as well.
#
# New tracking numbers are plain 12-digit numeric strings that must be
# written back as text (shared-string) cells, exactly like the existing
# values they replace - i.e. no visible number formatting / style change.
# Setting the NumberFormat to "@" (Text) forces the engine to store the
# value as a string instead of auto-coercing it to a number; resetting the
# cell's Style back to "Normal" afterwards drops the style index back to
# the default (unstyled) one the cell had before, so there's no residual
# formatting difference versus the original cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newTrackNums = @{
    2  = "320018628560"
    3  = "320018628571"
    4  = "320018628608"
    5  = "320018628620"
    6  = "320018628663"
    7  = "320018628685"
    8  = "320018628711"
    9  = "320018628733"
    10 = "320018628766"
    11 = "320018628788"
    12 = "320018628825"
    13 = "320018628847"
    14 = "320018628870"
    15 = "320018625274"
    16 = "320018625300"
    17 = "320018625322"
    18 = "320018625366"
    19 = "320018625388"
    20 = "320018625414"
    21 = "320018625436"
    22 = "320018625469"
}

# Rows where column D ("ShipmentTrackNum") mirrors column C ("PackageTrackNum")
$rowsWithD = @(5, 6, 7, 13, 14, 15, 16, 17)

foreach ($row in $newTrackNums.Keys) {
    $value = $newTrackNums[$row]

    $cellC = $ws.Cells.Item($row, 3)
    $cellC.NumberFormat = "@"
    $cellC.Value = $value
    $cellC.Style = "Normal"

    if ($rowsWithD -contains $row) {
        $cellD = $ws.Cells.Item($row, 4)
        $cellD.NumberFormat = "@"
        $cellD.Value = $value
        $cellD.Style = "Normal"
    }
}
